# Add 14 new players (rows 209-222) to the "Players Data - PSL 2025" sheet,
# continuing the existing REGISTRATION NO: / PLAYER NAME: / PLACE: / CATEGORY:
# / TEAMS PLAYED FOR table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPlayers = @(
    @(208, "KIRAN KAMAT",      "KUMTA",        "WK-BATSMAN",  "NONE"),
    @(209, "ADITYA BHAT",      "MANGALORE ",   "BOWLER",      "AVATAR ELEVEN MALPE"),
    @(210, "RAKSHITH SHENOY",  "MULKI",        "BOWLER",      "VOLALANKE FIGHTERS MULKI"),
    @(211, "VIVEK",            "AMBAGILU ",    "BATSMAN",     "UDUPI BLASTERS "),
    @(212, "BHUVAN BHAT",      "UDUPI",        "BATSMAN",     "JAIKAR STRIKERS"),
    @(213, "PANCHAM",          "MANGALORE",    "BATSMAN",     "SAPTHAMI WARRIORS"),
    @(214, "ADARSH BHAT ",     "UDUPI",        "ALL-ROUNDER", "VV WARRIORS "),
    @(215, "ASHLESH SHENOY",   "UDUPI",        "ALL-ROUNDER", "VEERANJANEYA CRICKETERS KAUP, SPARK MANGALAPURA,GPL 2023 UDUPI BLASTERS"),
    @(216, "ABHISHEK PAI",     "MULKI",        "ALL-ROUNDER", "DEADLY PANTHERS "),
    @(217, "ADESH ",           "MULKI ",       "BOWLER",      "VOLALANKE FIGHTERS "),
    @(218, "NITHIN KAMATH",    "MULKI",        "BATSMAN",     "VOLALANKE FIGHTERS MULKI"),
    @(219, "PRAJWAL KAMATH",   "UDUPI",        "ALL-ROUNDER", "AK KNIGHT RIDERS "),
    @(220, "VASANTH",          "UDUPI",        "ALL-ROUNDER", "UDUPI FRIENDS"),
    @(221, "VIKAS PAI",        "MANGALORE ",   "BATSMAN",     "TANK COBRA ")
)

$startRow = 209
for ($i = 0; $i -lt $newPlayers.Count; $i++) {
    $rowNum = $startRow + $i
    $data = $newPlayers[$i]

    $arr = New-Object 'object[,]' 1,5
    $arr[0,0] = $data[0]
    $arr[0,1] = $data[1]
    $arr[0,2] = $data[2]
    $arr[0,3] = $data[3]
    $arr[0,4] = $data[4]

    $ws.Range("A$rowNum" + ":E$rowNum").Value = $arr
    $ws.Rows.Item($rowNum).RowHeight = 56
}

# Restore the active selection to where the author left off editing.
[void]$ws.Range("B208").Select()
